$wb = $excel.ActiveWorkbook

# Turn off alerts so deleting a sheet doesn't prompt a confirmation dialog
$excel.DisplayAlerts = $false

# Remove the old "version 1" draft sheet entirely - the spec is now published
# as a single sheet (content of the former "version 2" sheet).
$wsOld = $wb.Worksheets.Item("version 1")
$wsOld.Delete()

$excel.DisplayAlerts = $true

# Rename the remaining "version 2" sheet to "Sheet1"
$ws = $wb.Worksheets.Item("version 2")
$ws.Name = "Sheet1"

# The long requirement text in row 2 now needs to wrap onto two lines like
# the other long rows, so re-fit its height
$ws.Rows.Item(2).AutoFit()

# Update the window selection to reflect the published sheet's cursor position
$ws.Range("H9").Select()
